$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3665.1738
$ws.Range("I64").Value = 3799.6667
$ws.Range("J64").Value = 3645
$ws.Range("K64").Value = 3799.6667
$ws.Range("L64").Value = 3645
$ws.Range("M64").Value = -3551.6667
$ws.Range("N64").Value = -4141

$ws.Range("H67").Value = 3665.1738
$ws.Range("I67").Value = 3799.6667
$ws.Range("J67").Value = 3645
$ws.Range("K67").Value = 3799.6667
$ws.Range("L67").Value = 3645
$ws.Range("M67").Value = -2941.6667
$ws.Range("N67").Value = -5361

$ws.Range("H132").Value = 1849.242
$ws.Range("I132").Value = 1387.5
$ws.Range("J132").Value = 3773.1667
$ws.Range("K132").Value = 4162.5
$ws.Range("L132").Value = 11319.5001
$ws.Range("M132").Value = -1632.5
$ws.Range("N132").Value = -16379.5001

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200

$ws.Range("H137").Value = 1458.4706
$ws.Range("I137").Value = 1389.5454
$ws.Range("J137").Value = 1584.8334
$ws.Range("K137").Value = 4168.6362
$ws.Range("L137").Value = 4754.5002
$ws.Range("M137").Value = -1618.6362
$ws.Range("N137").Value = -9854.5002

$ws.Range("H138").Value = 2145.805
$ws.Range("I138").Value = 1171.6964
$ws.Range("J138").Value = 4243.885
$ws.Range("K138").Value = 3515.0892
$ws.Range("L138").Value = 12731.655
$ws.Range("M138").Value = 1624.9108
$ws.Range("N138").Value = -23011.655

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15870.523
$ws.Range("I32").Value = 15558
$ws.Range("J32").Value = 19933.334
$ws.Range("K32").Value = 15558
$ws.Range("L32").Value = 19933.334
$ws.Range("M32").Value = -15271
$ws.Range("N32").Value = -20507.334

$ws.Range("H61").Value = 2463.4167
$ws.Range("I61").Value = 1747.1538
$ws.Range("J61").Value = 4325.7
$ws.Range("K61").Value = 1747.1538
$ws.Range("L61").Value = 4325.7
$ws.Range("M61").Value = -1535.1538
$ws.Range("N61").Value = -4749.7

$ws.Range("H63").Value = 3588.125
$ws.Range("I63").Value = 3529.2856
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 3529.2856
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -2843.2856
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 3588.125
$ws.Range("I66").Value = 3529.2856
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 17646.428
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -14214.428
$ws.Range("N66").Value = -26864

$ws.Range("H74").Value = 639.2439000000001
$ws.Range("I74").Value = 531.8158
$ws.Range("K74").Value = 531.8158
$ws.Range("M74").Value = 342.1842

$ws.Range("H77").Value = 639.2439000000001
$ws.Range("I77").Value = 531.8158
$ws.Range("K77").Value = 2659.079
$ws.Range("M77").Value = 1708.921

$ws.Range("H132").Value = 6100.514
$ws.Range("I132").Value = 7691.4287
$ws.Range("K132").Value = 23074.2861
$ws.Range("M132").Value = -20544.2861

$ws.Range("H136").Value = 2463.4167
$ws.Range("I136").Value = 1747.1538
$ws.Range("J136").Value = 4325.7
$ws.Range("K136").Value = 5241.4614
$ws.Range("L136").Value = 12977.1
$ws.Range("M136").Value = -2691.4614
$ws.Range("N136").Value = -18077.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 42085.715
$ws.Range("I12").Value = 25500
$ws.Range("J12").Value = 48720
$ws.Range("K12").Value = 25500
$ws.Range("L12").Value = 48720
$ws.Range("M12").Value = -25332
$ws.Range("N12").Value = -49056

$ws.Range("H86").Value = 254226.62
$ws.Range("I86").Value = 10003
$ws.Range("K86").Value = 10003
$ws.Range("M86").Value = -8880

$ws.Range("H89").Value = 254226.62
$ws.Range("I89").Value = 10003
$ws.Range("K89").Value = 50015
$ws.Range("M89").Value = -44399

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1737.0465
$ws.Range("I31").Value = 1052.7354
$ws.Range("J31").Value = 4322.222
$ws.Range("K31").Value = 1052.7354
$ws.Range("L31").Value = 4322.222
$ws.Range("M31").Value = -757.7354
$ws.Range("N31").Value = -4912.222

$ws.Range("H34").Value = 1737.0465
$ws.Range("I34").Value = 1052.7354
$ws.Range("J34").Value = 4322.222
$ws.Range("K34").Value = 1052.7354
$ws.Range("L34").Value = 4322.222
$ws.Range("M34").Value = -850.7354
$ws.Range("N34").Value = -4726.222

$ws.Range("H86").Value = 8740.75
$ws.Range("I86").Value = 13867.444
$ws.Range("J86").Value = 4546.1816
$ws.Range("K86").Value = 13867.444
$ws.Range("L86").Value = 4546.1816
$ws.Range("M86").Value = -12744.444
$ws.Range("N86").Value = -6792.1816

$ws.Range("H89").Value = 8740.75
$ws.Range("I89").Value = 13867.444
$ws.Range("J89").Value = 4546.1816
$ws.Range("K89").Value = 69337.22
$ws.Range("L89").Value = 22730.908
$ws.Range("M89").Value = -63721.22
$ws.Range("N89").Value = -33962.908

$ws.Range("H134").Value = 1674.8959
$ws.Range("I134").Value = 1050.4865
$ws.Range("K134").Value = 3151.4595
$ws.Range("M134").Value = -616.4594999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 783.86664
$ws.Range("I2").Value = 1522
$ws.Range("J2").Value = 138
$ws.Range("K2").Value = 9132
$ws.Range("L2").Value = 828
$ws.Range("M2").Value = -9019
$ws.Range("N2").Value = -1054

$ws.Range("H107").Value = 570
$ws.Range("J107").Value = 904
$ws.Range("L107").Value = 2712
$ws.Range("N107").Value = -6552

$ws.Range("H113").Value = 591.913
$ws.Range("I113").Value = 542.1667
$ws.Range("K113").Value = 1626.5001
$ws.Range("M113").Value = 543.4999

$ws.Range("H118").Value = 2596.2222
$ws.Range("I118").Value = 1067.5
$ws.Range("J118").Value = 3819.2
$ws.Range("K118").Value = 3202.5
$ws.Range("L118").Value = 11457.6
$ws.Range("M118").Value = -1959.5
$ws.Range("N118").Value = -13943.6

$ws.Range("H122").Value = 1189.8572
$ws.Range("I122").Value = 885
$ws.Range("K122").Value = 7965
$ws.Range("M122").Value = -5515

$ws.Range("H131").Value = 4518.763
$ws.Range("I131").Value = 8702.416999999999
$ws.Range("J131").Value = 2587.8462
$ws.Range("K131").Value = 26107.251
$ws.Range("L131").Value = 7763.5386
$ws.Range("M131").Value = -21067.251
$ws.Range("N131").Value = -17843.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3725.9092
$ws.Range("I80").Value = 4164.1665
$ws.Range("J80").Value = 3200
$ws.Range("K80").Value = 4164.1665
$ws.Range("L80").Value = 3200
$ws.Range("M80").Value = -3166.1665
$ws.Range("N80").Value = -5196

$ws.Range("H83").Value = 3725.9092
$ws.Range("I83").Value = 4164.1665
$ws.Range("J83").Value = 3200
$ws.Range("K83").Value = 20820.8325
$ws.Range("L83").Value = 16000
$ws.Range("M83").Value = -15828.8325
$ws.Range("N83").Value = -25984

$ws.Range("H116").Value = 49899
$ws.Range("J116").Value = 49899
$ws.Range("L116").Value = 49899
$ws.Range("N116").Value = -59077

$ws.Range("H135").Value = 49996.668
$ws.Range("J135").Value = 49996.668
$ws.Range("L135").Value = 49996.668
$ws.Range("N135").Value = -60136.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3442.4285
$ws.Range("I7").Value = 3816.5
$ws.Range("J7").Value = 3161.875
$ws.Range("K7").Value = 3816.5
$ws.Range("L7").Value = 3161.875
$ws.Range("M7").Value = -3704.5
$ws.Range("N7").Value = -3385.875

$ws.Range("H100").Value = 6378
$ws.Range("I100").Value = 9130
$ws.Range("K100").Value = 9130
$ws.Range("M100").Value = -8589

$ws.Range("H126").Value = 3442.4285
$ws.Range("I126").Value = 3816.5
$ws.Range("J126").Value = 3161.875
$ws.Range("K126").Value = 11449.5
$ws.Range("L126").Value = 9485.625
$ws.Range("M126").Value = -8979.5
$ws.Range("N126").Value = -14425.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1136.0488
$ws.Range("I136").Value = 1128.921
$ws.Range("J136").Value = 1226.3334
$ws.Range("K136").Value = 3386.763
$ws.Range("L136").Value = 3679.0002
$ws.Range("M136").Value = -836.7629999999999
$ws.Range("N136").Value = -8779.0002
